# Update the "Rules" worksheet: rename rule in B11 from "R40" to "1".
# A leading apostrophe forces Excel to store the value as text (shared
# string) rather than coercing the numeric-looking value to a number,
# matching the original cell's text type.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
